# Weekly refresh: a new price-report row is inserted at row 18, pushing the
# existing data (rows 18..85) down by one row (to rows 19..86).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 - shifts rows 18..85 down to 19..86.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new week's data.
$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(18, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(18, 3).Value = "La Araucanía"
$ws.Cells.Item(18, 4).Value = 44592
$ws.Cells.Item(18, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18, 5).Value = 9
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100101
$ws.Cells.Item(18, 8).Value = "Berries"
$ws.Cells.Item(18, 9).Value = 100101001
$ws.Cells.Item(18, 10).Value = "Arándano (blue)"
$ws.Cells.Item(18, 11).Value = "Sin especificar"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 20
$ws.Cells.Item(18, 14).Value = 1800
$ws.Cells.Item(18, 15).Value = 1800
$ws.Cells.Item(18, 16).Value = 1800
$ws.Cells.Item(18, 17).Value = "`$/kilo"
$ws.Cells.Item(18, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(18, 19).Value = 1800
$ws.Cells.Item(18, 20).Value = 1
